# Auto-generated script applying the Leviathan_Profits price-refresh diff.
# Updates currentAveragePrice / NQ / HQ derived columns (H-N) for the rows
# touched by the scheduled market-data runner, sheet by sheet.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 480.8
$ws.Range("I33").Value = 476
$ws.Range("K33").Value = 476
$ws.Range("M33").Value = -247
$ws.Range("H63").Value = 25271
$ws.Range("J63").Value = 25271
$ws.Range("L63").Value = 25271
$ws.Range("N63").Value = -26519
$ws.Range("H66").Value = 25271
$ws.Range("J66").Value = 25271
$ws.Range("L66").Value = 75813
$ws.Range("N66").Value = -82053
$ws.Range("H70").Value = 3618.7693
$ws.Range("I70").Value = 3416.75
$ws.Range("J70").Value = 3942
$ws.Range("K70").Value = 10250.25
$ws.Range("L70").Value = 11826
$ws.Range("M70").Value = -9980.25
$ws.Range("N70").Value = -12366
$ws.Range("H73").Value = 3618.7693
$ws.Range("I73").Value = 3416.75
$ws.Range("J73").Value = 3942
$ws.Range("K73").Value = 10250.25
$ws.Range("L73").Value = 11826
$ws.Range("M73").Value = -9314.25
$ws.Range("N73").Value = -13698
$ws.Range("H80").Value = 1751.8572
$ws.Range("I80").Value = 1014.8461
$ws.Range("K80").Value = 3044.5383
$ws.Range("M80").Value = -2046.5383
$ws.Range("H83").Value = 1751.8572
$ws.Range("I83").Value = 1014.8461
$ws.Range("K83").Value = 9133.6149
$ws.Range("M83").Value = -4141.6149
$ws.Range("H92").Value = 342.3
$ws.Range("I92").Value = 316.69232
$ws.Range("J92").Value = 389.85715
$ws.Range("K92").Value = 316.69232
$ws.Range("L92").Value = 389.85715
$ws.Range("M92").Value = 931.30768
$ws.Range("N92").Value = -2885.85715
$ws.Range("H137").Value = 1630.04
$ws.Range("I137").Value = 965.25
$ws.Range("K137").Value = 2895.75
$ws.Range("M137").Value = -345.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1594.098
$ws.Range("I61").Value = 1486.3556
$ws.Range("K61").Value = 1486.3556
$ws.Range("M61").Value = -1274.3556
$ws.Range("H74").Value = 2120
$ws.Range("I74").Value = 1732.2941
$ws.Range("K74").Value = 1732.2941
$ws.Range("M74").Value = -858.2941000000001
$ws.Range("H77").Value = 2120
$ws.Range("I77").Value = 1732.2941
$ws.Range("K77").Value = 8661.470499999999
$ws.Range("M77").Value = -4293.470499999999
$ws.Range("H132").Value = 7553.206
$ws.Range("I132").Value = 7832.0967
$ws.Range("K132").Value = 23496.2901
$ws.Range("M132").Value = -20966.2901
$ws.Range("H136").Value = 1594.098
$ws.Range("I136").Value = 1486.3556
$ws.Range("K136").Value = 4459.066800000001
$ws.Range("M136").Value = -1909.066800000001
$ws.Range("H140").Value = 110514.5
$ws.Range("J140").Value = 110514.5
$ws.Range("L140").Value = 110514.5
$ws.Range("N140").Value = -120874.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 31706.559
$ws.Range("I31").Value = 35423.465
$ws.Range("K31").Value = 35423.465
$ws.Range("M31").Value = -35128.465
$ws.Range("H34").Value = 31706.559
$ws.Range("I34").Value = 35423.465
$ws.Range("K34").Value = 35423.465
$ws.Range("M34").Value = -35221.465
$ws.Range("H58").Value = 1477.7059
$ws.Range("I58").Value = 1543.9
$ws.Range("J58").Value = 1383.1428
$ws.Range("K58").Value = 1543.9
$ws.Range("L58").Value = 1383.1428
$ws.Range("M58").Value = -1340.9
$ws.Range("N58").Value = -1789.1428
$ws.Range("H59").Value = 19999.6
$ws.Range("J59").Value = 19999.6
$ws.Range("L59").Value = 19999.6
$ws.Range("N59").Value = -22289.6
$ws.Range("H60").Value = 1000
$ws.Range("J60").Value = 0
$ws.Range("L60").Value = 0
$ws.Range("N60").ClearContents()
$ws.Range("H132").Value = 3722.2693
$ws.Range("I132").Value = 3804.3333
$ws.Range("J132").Value = 2737.5
$ws.Range("K132").Value = 11412.9999
$ws.Range("L132").Value = 8212.5
$ws.Range("M132").Value = -8882.999899999999
$ws.Range("N132").Value = -13272.5
$ws.Range("H134").Value = 2365.3618
$ws.Range("I134").Value = 1734.85
$ws.Range("K134").Value = 5204.549999999999
$ws.Range("M134").Value = -2669.549999999999
$ws.Range("H136").Value = 1477.7059
$ws.Range("I136").Value = 1543.9
$ws.Range("J136").Value = 1383.1428
$ws.Range("K136").Value = 4631.700000000001
$ws.Range("L136").Value = 4149.428400000001
$ws.Range("M136").Value = -2081.700000000001
$ws.Range("N136").Value = -9249.428400000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1459.7778
$ws.Range("J68").Value = 1441.3334
$ws.Range("L68").Value = 4324.0002
$ws.Range("N68").Value = -5946.0002
$ws.Range("H71").Value = 1459.7778
$ws.Range("J71").Value = 1441.3334
$ws.Range("L71").Value = 12972.0006
$ws.Range("N71").Value = -21084.0006
$ws.Range("H114").Value = 25001862
$ws.Range("I114").Value = 33335232
$ws.Range("J114").Value = 1750
$ws.Range("K114").Value = 100005696
$ws.Range("L114").Value = 5250
$ws.Range("M114").Value = -100002442
$ws.Range("N114").Value = -11758
$ws.Range("H117").Value = 5536.1
$ws.Range("J117").Value = 8400.25
$ws.Range("L117").Value = 25200.75
$ws.Range("N117").Value = -32084.75
$ws.Range("H122").Value = 500.2143
$ws.Range("I122").Value = 545.1111
$ws.Range("K122").Value = 4905.9999
$ws.Range("M122").Value = -2455.9999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 34000
$ws.Range("J52").Value = 34000
$ws.Range("L52").Value = 34000
$ws.Range("N52").Value = -34518
$ws.Range("H102").Value = 2506.1904
$ws.Range("I102").Value = 2506.9473
$ws.Range("J102").Value = 2499
$ws.Range("K102").Value = 2506.9473
$ws.Range("L102").Value = 2499
$ws.Range("M102").Value = -884.9472999999998
$ws.Range("N102").Value = -5743
$ws.Range("H107").Value = 19232668
$ws.Range("I107").Value = 768.63635
$ws.Range("K107").Value = 768.63635
$ws.Range("M107").Value = 1151.36365
$ws.Range("H126").Value = 2478.6
$ws.Range("J126").Value = 2222.8333
$ws.Range("L126").Value = 6668.499899999999
$ws.Range("N126").Value = -11608.4999
$ws.Range("H136").Value = 48735.734
$ws.Range("J136").Value = 48735.734
$ws.Range("L136").Value = 146207.202
$ws.Range("N136").Value = -151307.202

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 862.4737
$ws.Range("I22").Value = 623
$ws.Range("J22").Value = 973
$ws.Range("K22").Value = 623
$ws.Range("L22").Value = 973
$ws.Range("M22").Value = -328
$ws.Range("N22").Value = -1563
$ws.Range("H27").Value = 862.4737
$ws.Range("I27").Value = 623
$ws.Range("J27").Value = 973
$ws.Range("K27").Value = 623
$ws.Range("L27").Value = 973
$ws.Range("M27").Value = -516
$ws.Range("N27").Value = -1187
$ws.Range("H40").Value = 64864.65
$ws.Range("I40").Value = 26459.8
$ws.Range("K40").Value = 26459.8
$ws.Range("M40").Value = -26323.8
$ws.Range("H46").Value = 32381
$ws.Range("I46").Value = 54804.625
$ws.Range("J46").Value = 2482.8333
$ws.Range("K46").Value = 54804.625
$ws.Range("L46").Value = 2482.8333
$ws.Range("M46").Value = -54616.625
$ws.Range("N46").Value = -2858.8333
$ws.Range("H82").Value = 1979.5714
$ws.Range("I82").Value = 1605.6666
$ws.Range("K82").Value = 1605.6666
$ws.Range("M82").Value = -1244.6666
$ws.Range("H85").Value = 1979.5714
$ws.Range("I85").Value = 1605.6666
$ws.Range("K85").Value = 1605.6666
$ws.Range("M85").Value = -357.6666
$ws.Range("H122").Value = 158765.16
$ws.Range("I122").Value = 226068.44
$ws.Range("J122").Value = 7332.75
$ws.Range("K122").Value = 678205.3200000001
$ws.Range("L122").Value = 21998.25
$ws.Range("M122").Value = -675755.3200000001
$ws.Range("N122").Value = -26898.25
$ws.Range("H136").Value = 3133.389
$ws.Range("I136").Value = 2687.5334
$ws.Range("K136").Value = 8062.600199999999
$ws.Range("M136").Value = -5512.600199999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1913.5555
$ws.Range("I126").Value = 1913.5555
$ws.Range("K126").Value = 5740.666499999999
$ws.Range("M126").Value = -3270.666499999999
$ws.Range("H132").Value = 837922.7
$ws.Range("I132").Value = 2081.7715
$ws.Range("J132").Value = 2167669.5
$ws.Range("K132").Value = 6245.314499999999
$ws.Range("L132").Value = 6503008.5
$ws.Range("M132").Value = -3715.314499999999
$ws.Range("N132").Value = -6508068.5
$ws.Range("H136").Value = 1060.4445
$ws.Range("I136").Value = 1055.0769
$ws.Range("K136").Value = 3165.2307
$ws.Range("M136").Value = -615.2307000000001

